$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.794.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "'2.344.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'239.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "'0.668"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("D7").Value = "'72.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.94%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.09%  "
$ws.Range("D10").Value = "'0.0996"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").Value = "'58.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").Value = "'32.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  -4.65%  "
$ws.Range("D15").Value = "'2.690.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'16.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "'0.900"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "'2.346.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "'43.706.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "'78.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'253.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  +8.25%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'3.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.33%  "
$ws.Range("D27").Value = "'2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "'10.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("D29").Value = "'2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "'175.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").Value = "'5.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").Value = "'5.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'3.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").Value = "'5.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.07%  "
$ws.Range("D42").Value = "'64.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.43%  "
$ws.Range("D43").Value = "'9.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("E44").Value = "  +6.14%  "
$ws.Range("D45").Value = "'18.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "'2.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").Value = "'98.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("E51").Value = "  -5.06%  "
